# Workbook was re-uploaded/resaved in Excel:
#  - F1 header text changed from "syst_u(%)" to "syst_c(%)"
#  - active cell selection moved from G12 to J4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "syst_c(%)"

[void]$ws.Range("J4").Select()
